# This script applies the changes described by the diff to the "Artfynd" worksheet.
#
# Summary of the edit:
#   - Row 2 and Row 3 (a "Spillkråka" observation and an "Ullticka" observation) had their
#     content reordered/swapped (the specific columns that differ between them are swapped),
#     likely due to the source data being re-sorted.
#   - Row 22 and Row 23 (a "Sprickporing" observation and a "Garnlav" observation) were
#     similarly swapped.
#   - Every data row's "Taxonsorteringsordning" value (column B) increased by 4
#     (a taxonomic sort-order id bump coming from the source system), applied AFTER the
#     row swaps above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell {
    param($col, $row1, $row2)
    $c1 = $ws.Range("$col$row1")
    $c2 = $ws.Range("$col$row2")
    $v1 = $c1.Value2
    $v2 = $c2.Value2
    if ($v2 -eq $null) { $c1.ClearContents() } else { $c1.Value2 = $v2 }
    if ($v1 -eq $null) { $c2.ClearContents() } else { $c2.Value2 = $v1 }
}

# --- Swap the columns that differ between row 2 and row 3 ---
$cols2_3 = @("A","B","E","F","G","H","M","P","Q","R","Z","AB","AC")
foreach ($col in $cols2_3) {
    Swap-Cell $col 2 3
}

# --- Swap the columns that differ between row 22 and row 23 ---
$cols22_23 = @("A","B","D","E","F","G","H","P","Q","R","S","Z","AB","AC")
foreach ($col in $cols22_23) {
    Swap-Cell $col 22 23
}

# --- Apply the +4 bump to column B (Taxonsorteringsordning) for every data row ---
for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Range("B$r")
    $cell.Value2 = $cell.Value2 + 4
}
